$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - literal input values (reading the input, not formulas)
$ws.Range("A10").Value = "Test11"
$ws.Range("B10").Value = "Test111"
$ws.Range("C10").Value = 1111
$ws.Range("D10").Value = "TestCompany"
$ws.Range("E10").Value = "Test"
$ws.Range("F10").Value = 250005

# Row 11 - formulas referencing row 10 and other existing rows
$ws.Range("A11").Formula = "=A10"
$ws.Range("B11").Formula = "=B10"
$ws.Range("C11").Formula = "=C8"
$ws.Range("D11").Formula = "=D3"
$ws.Range("E11").Formula = "=E4"
$ws.Range("F11").Formula = "=F5"

# Row 12 - formulas referencing earlier rows (C12 stays a literal value)
$ws.Range("A12").Formula = "=A2"
$ws.Range("B12").Formula = "=B4"
$ws.Range("C12").Value = 23
$ws.Range("D12").Formula = "=D3"
$ws.Range("E12").Formula = "=E11"
$ws.Range("F12").Formula = "=F6"

# Update the active selection to match the author's final cursor position
$ws.Range("F13").Select()
